$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "63.077.92"
Set-TextValue "E2" "  -0.80%  "

Set-TextValue "D3" "3.055.64"
Set-TextValue "E3" "  -1.14%  "

Set-TextValue "E4" "  -0.16%  "

Set-TextValue "D5" "583.68"
Set-TextValue "E5" "  -1.41%  "

Set-TextValue "D6" "151.76"
Set-TextValue "E6" "  -2.36%  "

Set-TextValue "E7" "  -0.14%  "

Set-TextValue "D8" "0.535"
Set-TextValue "E8" "  -1.76%  "

Set-TextValue "D9" "3.057.07"
Set-TextValue "E9" "  -0.87%  "

Set-TextValue "D10" "0.154"
Set-TextValue "E10" "  -2.73%  "

Set-TextValue "E11" "  -0.03%  "

Set-TextValue "E12" "  -2.34%  "

Set-TextValue "D13" "0.0000235"
Set-TextValue "E13" "  -2.95%  "

Set-TextValue "D14" "36.17"
Set-TextValue "E14" "  -3.87%  "

Set-TextValue "D16" "3.554.52"
Set-TextValue "E16" "  -1.25%  "

Set-TextValue "D17" "7.16"
Set-TextValue "E17" "  -0.87%  "

Set-TextValue "D18" "63.020.14"
Set-TextValue "E18" "  -0.84%  "

Set-TextValue "D19" "3.053.99"
Set-TextValue "E19" "  -0.94%  "

Set-TextValue "D20" "482.20"
Set-TextValue "E20" "  +1.19%  "

Set-TextValue "D21" "14.32"
Set-TextValue "E21" "  -2.77%  "

Set-TextValue "D22" "0.708"
Set-TextValue "E22" "  -1.64%  "

Set-TextValue "D23" "7.52"
Set-TextValue "E23" "  -0.66%  "

Set-TextValue "D24" "2.40"
Set-TextValue "E24" "  -0.81%  "

Set-TextValue "D25" "82.02"
Set-TextValue "E25" "  +0.91%  "

Set-TextValue "D26" "12.67"
Set-TextValue "E26" "  -2.23%  "

Set-TextValue "D27" "10.58"
Set-TextValue "E27" "  +4.79%  "

Set-TextValue "E28" "  +0.02%  "

Set-TextValue "D29" "7.40"
Set-TextValue "E29" "  +0.30%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "2.66"
Set-TextValue "E30" "  -1.41%  "

$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D31" "1.00"
Set-TextValue "E31" "  -0.15%  "

Set-TextValue "E32" "  +0.35%  "

Set-TextValue "D33" "27.81"
Set-TextValue "E33" "  +1.89%  "

Set-TextValue "E34" "  -2.70%  "

Set-TextValue "E35" "  +0.80%  "

Set-TextValue "E36" "  -4.31%  "

Set-TextValue "D37" "5.92"
Set-TextValue "E37" "  -3.44%  "

Set-TextValue "D38" "3.23"
Set-TextValue "E38" "  -5.13%  "

Set-TextValue "E39" "  -0.93%  "

Set-TextValue "D40" "9.23"
Set-TextValue "E40" "  -1.56%  "

Set-TextValue "D41" "50.44"
Set-TextValue "E41" "  -0.71%  "

Set-TextValue "D42" "428.63"
Set-TextValue "E42" "  -3.58%  "

Set-TextValue "D43" "0.287"
Set-TextValue "E43" "  +0.50%  "

Set-TextValue "E44" "  +3.75%  "

Set-TextValue "D46" "2.843.26"
Set-TextValue "E46" "  +1.35%  "

Set-TextValue "D47" "38.16"
Set-TextValue "E47" "  -4.97%  "

Set-TextValue "D48" "127.29"
Set-TextValue "E48" "  -3.54%  "

Set-TextValue "E49" "  +0.00%  "

Set-TextValue "D50" "25.10"
Set-TextValue "E50" "  -1.47%  "

Set-TextValue "D51" "0.110"
Set-TextValue "E51" "  -1.25%  "
